# Inserta un nuevo registro semanal de precios (Papa, Macroferia Regional de
# Talca) en la fila 386, desplazando hacia abajo las filas existentes
# 386-412 (que pasan a ser 387-413), tal como lo describe el diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserta una fila nueva en la posición 386; Excel copia el formato de la
# fila de arriba (incluye el estilo de fecha de la columna D) y desplaza
# todo lo de abajo (incluida la fila 412) una posición hacia abajo.
$ws.Rows(386).Insert()

$row = 386
$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 44610
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = 100114001
$ws.Cells.Item($row, 7).Value  = "Papa"
$ws.Cells.Item($row, 8).Value  = "Patagonia"
$ws.Cells.Item($row, 9).Value  = "1a nueva(o)"
$ws.Cells.Item($row, 10).Value = 1300
$ws.Cells.Item($row, 11).Value = 6000
$ws.Cells.Item($row, 12).Value = 6000
$ws.Cells.Item($row, 13).Value = 6000
$ws.Cells.Item($row, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($row, 15).Value = "Región del Maule"
$ws.Cells.Item($row, 16).Value = 240
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
